$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (QUILT, second occurrence): Progress 50 -> 62.5
$ws.Range("C4").Value = 62.5

# Row 9 (ALLEGRETTO-LTE (B7981028)): Days remaining 3 -> 2
$ws.Range("B9").Value = 2

# Row 11 (REJOICE (MK-5909-003)): Days remaining 30 -> 29
$ws.Range("B11").Value = 29

$wb.Save()
